$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("hotel_info")

$ws1.Range("G2:I2").NumberFormat = "@"
$ws1.Range("G2").Value = "5"
$ws1.Range("H2").Value = "462"
$ws1.Range("I2").Value = "5"
$ws1.Range("G2:I2").Style = "Normal"
